# Update Name of Algo
# Applies updated RandomForest imputation results to specific cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = -7.370799999999996
$ws.Range("A3").Value = -22.0431
$ws.Range("A14").Value = -21.85200000000001
$ws.Range("A21").Value = -20.02459999999999
$ws.Range("A23").Value = -20.63999999999998
$ws.Range("A25").Value = -21.82059999999999
$ws.Range("D25").Value = -8.301600000000002
$ws.Range("A26").Value = -21.01799999999996
$ws.Range("D27").Value = -8.835500000000007
$ws.Range("A29").Value = -20.88259999999998
$ws.Range("D31").Value = -8.497800000000007
$ws.Range("D39").Value = -7.850199999999997
$ws.Range("D48").Value = -7.337899999999998
$ws.Range("D51").Value = -7.818099999999999
$ws.Range("D52").Value = -7.721899999999997
$ws.Range("A53").Value = -21.83379999999999
$ws.Range("D55").Value = -8.104599999999994
$ws.Range("D56").Value = -7.892599999999998
$ws.Range("A57").Value = -22.35990000000002
$ws.Range("D57").Value = -8.390700000000002
$ws.Range("A59").Value = -22.2342
$ws.Range("A69").Value = -21.64399999999999
$ws.Range("D73").Value = -7.934499999999996
$ws.Range("A79").Value = -20.17170000000002
$ws.Range("A83").Value = -21.91329999999999
$ws.Range("D89").Value = -6.014100000000002
$ws.Range("D90").Value = -8.097800000000001
$ws.Range("A91").Value = -21.38030000000002
$ws.Range("D92").Value = -6.161800000000002
$ws.Range("A93").Value = -21.01779999999998
